$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D:E columns so numeric-looking strings (e.g. "0.9994")
# are not auto-converted to numbers by Excel, matching the original inlineStr text cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.274.55'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '1.839.77'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '233.43'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').Value = '0.9995'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.4645'
$ws.Range('E7').Value = '  -2.73%  '
$ws.Range('D8').Value = '0.2726'
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('D9').Value = '0.06273'
$ws.Range('E9').Value = '  -3.71%  '
$ws.Range('D10').Value = '1.837.52'
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').Value = '0.07415'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').Value = '16.27'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').Value = '4.931'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('D14').Value = '83.61'
$ws.Range('E14').Value = '  -3.58%  '
$ws.Range('D15').Value = '0.6201'
$ws.Range('E15').Value = '  -2.94%  '
$ws.Range('D16').Value = '30.195.52'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '0.9999'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '227.83'
$ws.Range('E18').Value = '  -2.66%  '
$ws.Range('D19').Value = '0.000007283'
$ws.Range('E19').Value = '  -2.11%  '
$ws.Range('D20').Value = '12.32'
$ws.Range('E20').Value = '  -4.81%  '
$ws.Range('D21').Value = '0.9997'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '4.894'
$ws.Range('E22').Value = '  -4.32%  '
$ws.Range('D23').Value = '5.864'
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('D24').Value = '9.168'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('D25').Value = '163.99'
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('D26').Value = '17.78'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').Value = '1.867'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('D28').Value = '0.1035'
$ws.Range('E28').Value = '  -1.35%  '
$ws.Range('D29').Value = '1.368'
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('D30').Value = '4.074'
$ws.Range('E30').Value = '  -4.38%  '
$ws.Range('D31').Value = '3.801'
$ws.Range('E31').Value = '  -4.02%  '
$ws.Range('D32').Value = '0.04838'
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('D33').Value = '1.138'
$ws.Range('E33').Value = '  -2.36%  '
$ws.Range('D34').Value = '0.7058'
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('D35').Value = '2.689'
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').Value = '0.01879'
$ws.Range('E36').Value = '  -3.01%  '
$ws.Range('D37').Value = '2.648'
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('D38').Value = '0.8859'
$ws.Range('E38').Value = '  -2.79%  '
$ws.Range('D39').Value = '104.81'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').Value = '1.918'
$ws.Range('E40').Value = '  -5.34%  '
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('D42').Value = '5.549'
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').Value = '0.4011'
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('D44').Value = '7.043'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '60.19'
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '0.1196'
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('D47').Value = '8.582'
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('D48').Value = '33.11'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').Value = '0.05507'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('D50').Value = '1.351'
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('D51').Value = '0.3634'
$ws.Range('E51').Value = '  -3.32%  '

# Restore default (Normal) style so no extra number-format styling is introduced
$ws.Range("D2:E51").Style = "Normal"
